$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 878.4286
$ws.Range("I2").Value = 449.66666
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 449.66666
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = -336.66666
$ws.Range("N2").Value = -1426
$ws.Range("H28").Value = 11899.8
$ws.Range("I28").Value = 1098.75
$ws.Range("K28").Value = 1098.75
$ws.Range("M28").Value = -613.75
$ws.Range("H38").Value = 63.5
$ws.Range("I38").Value = 63.5
$ws.Range("K38").Value = 190.5
$ws.Range("M38").Value = 181.5
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 9000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -10872
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 7537.8
$ws.Range("I74").Value = 7172.25
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 7172.25
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -6236.25
$ws.Range("N74").Value = -10872
$ws.Range("H77").Value = 7537.8
$ws.Range("I77").Value = 7172.25
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 35861.25
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -31181.25
$ws.Range("N77").Value = -54360
$ws.Range("H80").Value = 620
$ws.Range("I80").Value = 541.6667
$ws.Range("K80").Value = 1625.0001
$ws.Range("M80").Value = -627.0001
$ws.Range("H83").Value = 620
$ws.Range("I83").Value = 541.6667
$ws.Range("K83").Value = 4875.0003
$ws.Range("M83").Value = 116.9997000000003
$ws.Range("H94").Value = 4054.111
$ws.Range("I94").Value = 4054.111
$ws.Range("K94").Value = 4054.111
$ws.Range("M94").Value = -3603.111
$ws.Range("H111").Value = 1759.4073
$ws.Range("I111").Value = 627
$ws.Range("J111").Value = 1850
$ws.Range("K111").Value = 1881
$ws.Range("L111").Value = 5550
$ws.Range("M111").Value = 1186
$ws.Range("N111").Value = -11684
$ws.Range("H112").Value = 3410.7778
$ws.Range("J112").Value = 3385.4285
$ws.Range("L112").Value = 10156.2855
$ws.Range("N112").Value = -12372.2855
$ws.Range("H129").Value = 970.4666999999999
$ws.Range("I129").Value = 715.4167
$ws.Range("J129").Value = 1990.6666
$ws.Range("K129").Value = 2146.2501
$ws.Range("L129").Value = 5971.9998
$ws.Range("M129").Value = 2853.7499
$ws.Range("N129").Value = -15971.9998
$ws.Range("H138").Value = 1968.409
$ws.Range("J138").Value = 7629.25
$ws.Range("L138").Value = 22887.75
$ws.Range("N138").Value = -33167.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 2750
$ws.Range("I46").Value = 2750
$ws.Range("K46").Value = 2750
$ws.Range("M46").Value = -2431
$ws.Range("H61").Value = 1305.0714
$ws.Range("I61").Value = 1052
$ws.Range("K61").Value = 1052
$ws.Range("M61").Value = -840
$ws.Range("H97").Value = 657.2
$ws.Range("I97").Value = 717.5714
$ws.Range("J97").Value = 516.3333
$ws.Range("K97").Value = 717.5714
$ws.Range("L97").Value = 516.3333
$ws.Range("M97").Value = -221.5714
$ws.Range("N97").Value = -1508.3333
$ws.Range("H109").Value = 80000
$ws.Range("J109").Value = 80000
$ws.Range("L109").Value = 80000
$ws.Range("N109").Value = -82774
$ws.Range("H132").Value = 2754.5334
$ws.Range("I132").Value = 2729.2222
$ws.Range("K132").Value = 8187.6666
$ws.Range("M132").Value = -5657.6666
$ws.Range("H136").Value = 1305.0714
$ws.Range("I136").Value = 1052
$ws.Range("K136").Value = 3156
$ws.Range("M136").Value = -606
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 38490
$ws.Range("I26").Value = 38490
$ws.Range("K26").Value = 38490
$ws.Range("M26").Value = -38198
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H96").Value = 20212.5
$ws.Range("I96").Value = 20212.5
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 20212.5
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -17466.5
$ws.Range("N96").ClearContents()
$ws.Range("H105").Value = 2152
$ws.Range("I105").Value = 1815
$ws.Range("K105").Value = 1815
$ws.Range("M105").Value = -68
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 5887.1
$ws.Range("I6").Value = 181.28572
$ws.Range("J6").Value = 19200.666
$ws.Range("K6").Value = 181.28572
$ws.Range("L6").Value = 19200.666
$ws.Range("M6").Value = -68.28572
$ws.Range("N6").Value = -19426.666
$ws.Range("H7").Value = 103.73684
$ws.Range("I7").Value = 43.375
$ws.Range("J7").Value = 147.63637
$ws.Range("K7").Value = 43.375
$ws.Range("L7").Value = 147.63637
$ws.Range("M7").Value = 69.625
$ws.Range("N7").Value = -373.63637
$ws.Range("H16").Value = 1099.6666
$ws.Range("I16").Value = 1099.6666
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1099.6666
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -812.6666
$ws.Range("N16").ClearContents()
$ws.Range("H17").Value = 3283.0908
$ws.Range("I17").Value = 2679.3333
$ws.Range("K17").Value = 2679.3333
$ws.Range("M17").Value = -2505.3333
$ws.Range("H22").Value = 1100.2858
$ws.Range("I22").Value = 993.1667
$ws.Range("J22").Value = 1180.625
$ws.Range("K22").Value = 993.1667
$ws.Range("L22").Value = 1180.625
$ws.Range("M22").Value = -643.1667
$ws.Range("N22").Value = -1880.625
$ws.Range("H25").Value = 2937.1428
$ws.Range("I25").Value = 2760
$ws.Range("K25").Value = 2760
$ws.Range("M25").Value = -2586
$ws.Range("H58").Value = 3267.7856
$ws.Range("I58").Value = 2313.25
$ws.Range("K58").Value = 2313.25
$ws.Range("M58").Value = -2110.25
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H88").Value = 10500
$ws.Range("J88").Value = 10500
$ws.Range("L88").Value = 10500
$ws.Range("N88").Value = -11312
$ws.Range("H91").Value = 10500
$ws.Range("J91").Value = 10500
$ws.Range("L91").Value = 10500
$ws.Range("N91").Value = -13308
$ws.Range("H105").Value = 1510.6666
$ws.Range("I105").Value = 791.25
$ws.Range("J105").Value = 2949.5
$ws.Range("K105").Value = 791.25
$ws.Range("L105").Value = 2949.5
$ws.Range("M105").Value = 955.75
$ws.Range("N105").Value = -6443.5
$ws.Range("H113").Value = 1099.6666
$ws.Range("I113").Value = 1099.6666
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1099.6666
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1070.3334
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 2715.158
$ws.Range("J132").Value = 3760
$ws.Range("L132").Value = 11280
$ws.Range("N132").Value = -16340
$ws.Range("H136").Value = 3267.7856
$ws.Range("I136").Value = 2313.25
$ws.Range("K136").Value = 6939.75
$ws.Range("M136").Value = -4389.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 51.77778
$ws.Range("I2").Value = 45.714287
$ws.Range("J2").Value = 55.636364
$ws.Range("K2").Value = 274.285722
$ws.Range("L2").Value = 333.818184
$ws.Range("M2").Value = -161.285722
$ws.Range("N2").Value = -559.818184
$ws.Range("H7").Value = 21.083334
$ws.Range("I7").Value = 19.8
$ws.Range("J7").Value = 27.5
$ws.Range("K7").Value = 59.40000000000001
$ws.Range("L7").Value = 82.5
$ws.Range("M7").Value = 52.59999999999999
$ws.Range("N7").Value = -306.5
$ws.Range("H12").Value = 140.25
$ws.Range("J12").Value = 149.45454
$ws.Range("L12").Value = 448.36362
$ws.Range("N12").Value = -794.3636200000001
$ws.Range("H15").Value = 103.833336
$ws.Range("I15").Value = 87
$ws.Range("J15").Value = 112.25
$ws.Range("K15").Value = 261
$ws.Range("L15").Value = 336.75
$ws.Range("M15").Value = -121
$ws.Range("N15").Value = -616.75
$ws.Range("H57").Value = 1721.4615
$ws.Range("I57").Value = 1614.9166
$ws.Range("K57").Value = 4844.7498
$ws.Range("M57").Value = -4285.7498
$ws.Range("H68").Value = 701
$ws.Range("J68").Value = 701.5
$ws.Range("L68").Value = 2104.5
$ws.Range("N68").Value = -3726.5
$ws.Range("H71").Value = 701
$ws.Range("J71").Value = 701.5
$ws.Range("L71").Value = 6313.5
$ws.Range("N71").Value = -14425.5
$ws.Range("H98").Value = 311.1111
$ws.Range("I98").Value = 266.66666
$ws.Range("J98").Value = 333.33334
$ws.Range("K98").Value = 799.9999799999999
$ws.Range("L98").Value = 1000.00002
$ws.Range("M98").Value = 698.0000200000001
$ws.Range("N98").Value = -3996.00002
$ws.Range("H103").Value = 1672.9286
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 1672.9286
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 5018.7858
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -6776.7858
$ws.Range("H113").Value = 1634.5
$ws.Range("I113").Value = 806.6667
$ws.Range("K113").Value = 2420.0001
$ws.Range("M113").Value = -250.0001000000002
$ws.Range("H114").Value = 2795.4285
$ws.Range("J114").Value = 3124.3333
$ws.Range("L114").Value = 9372.999899999999
$ws.Range("N114").Value = -15880.9999
$ws.Range("H116").Value = 1924.1666
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 1386.25
$ws.Range("K116").Value = 9000
$ws.Range("L116").Value = 4158.75
$ws.Range("M116").Value = -5558
$ws.Range("N116").Value = -11042.75
$ws.Range("H122").Value = 485.125
$ws.Range("I122").Value = 373.93332
$ws.Range("K122").Value = 3365.39988
$ws.Range("M122").Value = -915.3998799999999
$ws.Range("H125").Value = 300
$ws.Range("I125").Value = 300
$ws.Range("K125").Value = 900
$ws.Range("M125").Value = 4020
$ws.Range("H129").Value = 2938.6667
$ws.Range("I129").Value = 1149.75
$ws.Range("J129").Value = 6516.5
$ws.Range("K129").Value = 3449.25
$ws.Range("L129").Value = 19549.5
$ws.Range("M129").Value = 1550.75
$ws.Range("N129").Value = -29549.5
$ws.Range("H138").Value = 4989.857
$ws.Range("J138").Value = 8266.666999999999
$ws.Range("L138").Value = 24800.001
$ws.Range("N138").Value = -35080.001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 724.8
$ws.Range("I97").Value = 724.8
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 724.8
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -228.8
$ws.Range("N97").ClearContents()
$ws.Range("H101").Value = 17164.334
$ws.Range("J101").Value = 17164.334
$ws.Range("L101").Value = 17164.334
$ws.Range("N101").Value = -23654.334
$ws.Range("H104").Value = 4341.5713
$ws.Range("J104").Value = 4341.5713
$ws.Range("L104").Value = 4341.5713
$ws.Range("N104").Value = -11329.5713
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2469.6667
$ws.Range("I22").Value = 701
$ws.Range("J22").Value = 2690.75
$ws.Range("K22").Value = 701
$ws.Range("L22").Value = 2690.75
$ws.Range("M22").Value = -406
$ws.Range("N22").Value = -3280.75
$ws.Range("H27").Value = 2469.6667
$ws.Range("I27").Value = 701
$ws.Range("J27").Value = 2690.75
$ws.Range("K27").Value = 701
$ws.Range("L27").Value = 2690.75
$ws.Range("M27").Value = -594
$ws.Range("N27").Value = -2904.75
$ws.Range("H46").Value = 2235
$ws.Range("I46").Value = 432.5
$ws.Range("J46").Value = 3136.25
$ws.Range("K46").Value = 432.5
$ws.Range("L46").Value = 3136.25
$ws.Range("M46").Value = -244.5
$ws.Range("N46").Value = -3512.25
$ws.Range("H61").Value = 2706.2
$ws.Range("I61").Value = 1872.091
$ws.Range("K61").Value = 1872.091
$ws.Range("M61").Value = -1670.091
$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248
$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240
$ws.Range("H82").Value = 2583
$ws.Range("I82").Value = 2583
$ws.Range("K82").Value = 2583
$ws.Range("M82").Value = -2222
$ws.Range("H85").Value = 2583
$ws.Range("I85").Value = 2583
$ws.Range("K85").Value = 2583
$ws.Range("M85").Value = -1335
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("H113").Value = 2706.2
$ws.Range("I113").Value = 1872.091
$ws.Range("K113").Value = 1872.091
$ws.Range("M113").Value = 297.9090000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 4000
$ws.Range("J97").Value = 4000
$ws.Range("L97").Value = 4000
$ws.Range("N97").Value = -5982
$ws.Range("H113").Value = 1439.8
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 1424.75
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 4274.25
$ws.Range("M113").Value = -2330
$ws.Range("N113").Value = -8614.25
$ws.Range("H132").Value = 3628.4285
$ws.Range("I132").Value = 3725
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 11175
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -8645
$ws.Range("N132").Value = -15559.0001
